# Auto-generated Excel COM-interop script
# Applies updated market-price values (currentAveragePrice / Leve profit columns)
# to the Sheets workbook, as produced by the scheduled pricing runner.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 12.5
$ws.Range("I11").Value = 12.5
$ws.Range("K11").Value = 12.5
$ws.Range("M11").Value = 127.5
$ws.Range("H33").Value = 306.8421
$ws.Range("I33").Value = 307.93332
$ws.Range("K33").Value = 307.93332
$ws.Range("M33").Value = -78.93331999999998
$ws.Range("H40").Value = 3473.762
$ws.Range("J40").Value = 5644.9
$ws.Range("L40").Value = 5644.9
$ws.Range("N40").Value = -5994.9
$ws.Range("H76").Value = 3666.3333
$ws.Range("I76").Value = 3666.3333
$ws.Range("K76").Value = 3666.3333
$ws.Range("M76").Value = -3351.3333
$ws.Range("H79").Value = 3666.3333
$ws.Range("I79").Value = 3666.3333
$ws.Range("K79").Value = 3666.3333
$ws.Range("M79").Value = -2574.3333
$ws.Range("H107").Value = 1358.5
$ws.Range("I107").Value = 1059.5
$ws.Range("J107").Value = 1508
$ws.Range("K107").Value = 1059.5
$ws.Range("L107").Value = 1508
$ws.Range("M107").Value = 860.5
$ws.Range("N107").Value = -5348
$ws.Range("H116").Value = 4364.4165
$ws.Range("J116").Value = 3624
$ws.Range("L116").Value = 3624
$ws.Range("N116").Value = -10508
$ws.Range("H118").Value = 1096.3334
$ws.Range("I118").Value = 1096.3334
$ws.Range("K118").Value = 3289.0002
$ws.Range("M118").Value = -1632.0002
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("M125").ClearContents()
$ws.Range("N125").ClearContents()
$ws.Range("H129").Value = 2689.2307
$ws.Range("I129").Value = 1827.5
$ws.Range("J129").Value = 3427.8572
$ws.Range("K129").Value = 5482.5
$ws.Range("L129").Value = 10283.5716
$ws.Range("M129").Value = -482.5
$ws.Range("N129").Value = -20283.5716
$ws.Range("H138").Value = 1471.3
$ws.Range("I138").Value = 1448.9231
$ws.Range("K138").Value = 4346.7693
$ws.Range("M138").Value = 793.2307000000001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11149.706
$ws.Range("I32").Value = 10303.066
$ws.Range("K32").Value = 10303.066
$ws.Range("M32").Value = -10016.066
$ws.Range("H45").Value = 2798.5
$ws.Range("I45").Value = 2731.6667
$ws.Range("K45").Value = 2731.6667
$ws.Range("M45").Value = -2354.6667
$ws.Range("H53").Value = 9500
$ws.Range("I53").Value = 9500
$ws.Range("K53").Value = 9500
$ws.Range("M53").Value = -8818
$ws.Range("H63").Value = 9249.833000000001
$ws.Range("I63").Value = 8375
$ws.Range("J63").Value = 10999.5
$ws.Range("K63").Value = 8375
$ws.Range("L63").Value = 10999.5
$ws.Range("M63").Value = -7689
$ws.Range("N63").Value = -12371.5
$ws.Range("H66").Value = 9249.833000000001
$ws.Range("I66").Value = 8375
$ws.Range("J66").Value = 10999.5
$ws.Range("K66").Value = 41875
$ws.Range("L66").Value = 54997.5
$ws.Range("M66").Value = -38443
$ws.Range("N66").Value = -61861.5
$ws.Range("H88").Value = 1753.25
$ws.Range("I88").Value = 1321.2
$ws.Range("K88").Value = 1321.2
$ws.Range("M88").Value = -915.2
$ws.Range("H91").Value = 1753.25
$ws.Range("I91").Value = 1321.2
$ws.Range("K91").Value = 1321.2
$ws.Range("M91").Value = 82.79999999999995
$ws.Range("H132").Value = 2811.2942
$ws.Range("I132").Value = 2811.2942
$ws.Range("K132").Value = 8433.882599999999
$ws.Range("M132").Value = -5903.882599999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2590.6365
$ws.Range("I105").Value = 2510.7778
$ws.Range("K105").Value = 2510.7778
$ws.Range("M105").Value = -763.7777999999998
$ws.Range("H134").Value = 1561.7
$ws.Range("I134").Value = 1561.7
$ws.Range("K134").Value = 4685.1
$ws.Range("M134").Value = -2150.1

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 492.22223
$ws.Range("I7").Value = 169.5
$ws.Range("J7").Value = 750.4
$ws.Range("K7").Value = 169.5
$ws.Range("L7").Value = 750.4
$ws.Range("M7").Value = -56.5
$ws.Range("N7").Value = -976.4
$ws.Range("H22").Value = 2354828.5
$ws.Range("I22").Value = 1813.4
$ws.Range("K22").Value = 1813.4
$ws.Range("M22").Value = -1463.4
$ws.Range("H35").Value = 1666.4445
$ws.Range("I35").Value = 1249.75
$ws.Range("J35").Value = 5000
$ws.Range("K35").Value = 1249.75
$ws.Range("L35").Value = 5000
$ws.Range("M35").Value = -955.75
$ws.Range("N35").Value = -5588
$ws.Range("H62").Value = 2647.6667
$ws.Range("J62").Value = 2448.5
$ws.Range("L62").Value = 2448.5
$ws.Range("N62").Value = -3696.5
$ws.Range("H65").Value = 2647.6667
$ws.Range("J65").Value = 2448.5
$ws.Range("L65").Value = 12242.5
$ws.Range("N65").Value = -18482.5
$ws.Range("H105").Value = 1855.8
$ws.Range("I105").Value = 1399.5
$ws.Range("K105").Value = 1399.5
$ws.Range("M105").Value = 347.5
$ws.Range("H141").Value = 393261.66
$ws.Range("J141").Value = 544892.5
$ws.Range("L141").Value = 544892.5
$ws.Range("N141").Value = -555252.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 71.888885
$ws.Range("I2").Value = 80.07143000000001
$ws.Range("J2").Value = 43.25
$ws.Range("K2").Value = 480.42858
$ws.Range("L2").Value = 259.5
$ws.Range("M2").Value = -367.42858
$ws.Range("N2").Value = -485.5
$ws.Range("H7").Value = 85972.75
$ws.Range("I7").Value = 250021.25
$ws.Range("J7").Value = 3948.5
$ws.Range("K7").Value = 750063.75
$ws.Range("L7").Value = 11845.5
$ws.Range("M7").Value = -749951.75
$ws.Range("N7").Value = -12069.5
$ws.Range("H107").Value = 153.75
$ws.Range("J107").Value = 153.75
$ws.Range("L107").Value = 461.25
$ws.Range("N107").Value = -4301.25
$ws.Range("H109").Value = 3430
$ws.Range("I109").Value = 4500
$ws.Range("K109").Value = 13500
$ws.Range("M109").Value = -12460

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 12666.444
$ws.Range("I43").Value = 916.6667
$ws.Range("J43").Value = 36166
$ws.Range("K43").Value = 916.6667
$ws.Range("L43").Value = 36166
$ws.Range("M43").Value = -765.6667
$ws.Range("N43").Value = -36468
$ws.Range("H70").Value = 8801.6
$ws.Range("I70").Value = 8754
$ws.Range("J70").Value = 8833.333000000001
$ws.Range("K70").Value = 8754
$ws.Range("L70").Value = 8833.333000000001
$ws.Range("M70").Value = -8484
$ws.Range("N70").Value = -9373.333000000001
$ws.Range("H73").Value = 8801.6
$ws.Range("I73").Value = 8754
$ws.Range("J73").Value = 8833.333000000001
$ws.Range("K73").Value = 8754
$ws.Range("L73").Value = 8833.333000000001
$ws.Range("M73").Value = -7818
$ws.Range("N73").Value = -10705.333
$ws.Range("H80").Value = 2937.4
$ws.Range("I80").Value = 3071.75
$ws.Range("K80").Value = 3071.75
$ws.Range("M80").Value = -2073.75
$ws.Range("H83").Value = 2937.4
$ws.Range("I83").Value = 3071.75
$ws.Range("K83").Value = 15358.75
$ws.Range("M83").Value = -10366.75
$ws.Range("H97").Value = 1730
$ws.Range("I97").Value = 1650
$ws.Range("K97").Value = 1650
$ws.Range("M97").Value = -1154
$ws.Range("H132").Value = 2881.111
$ws.Range("I132").Value = 2890.5
$ws.Range("J132").Value = 2848.25
$ws.Range("K132").Value = 8671.5
$ws.Range("L132").Value = 8544.75
$ws.Range("M132").Value = -6141.5
$ws.Range("N132").Value = -13604.75
$ws.Range("H134").Value = 34333
$ws.Range("J134").Value = 34333
$ws.Range("L134").Value = 102999
$ws.Range("N134").Value = -108069
$ws.Range("H136").Value = 65628.664
$ws.Range("J136").Value = 65628.664
$ws.Range("L136").Value = 196885.992
$ws.Range("N136").Value = -201985.992

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2950.3845
$ws.Range("I40").Value = 2335.6
$ws.Range("K40").Value = 2335.6
$ws.Range("M40").Value = -2199.6
$ws.Range("H55").Value = 515.4545000000001
$ws.Range("J55").Value = 1739.5
$ws.Range("L55").Value = 1739.5
$ws.Range("N55").Value = -2085.5
$ws.Range("H132").Value = 4961.25
$ws.Range("I132").Value = 4955.7144
$ws.Range("K132").Value = 14867.1432
$ws.Range("M132").Value = -12337.1432

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("M15").ClearContents()
$ws.Range("H81").Value = 1854.6
$ws.Range("I81").Value = 1977.6666
$ws.Range("K81").Value = 3955.3332
$ws.Range("M81").Value = -2894.3332
$ws.Range("H84").Value = 1854.6
$ws.Range("I84").Value = 1977.6666
$ws.Range("K84").Value = 19776.666
$ws.Range("M84").Value = -14472.666
$ws.Range("H113").Value = 424.08334
$ws.Range("I113").Value = 222.25
$ws.Range("K113").Value = 666.75
$ws.Range("M113").Value = 1503.25
$ws.Range("H132").Value = 2660.375
$ws.Range("I132").Value = 2220.4092
$ws.Range("K132").Value = 6661.2276
$ws.Range("M132").Value = -4131.2276
